$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.416622757911682
$ws.Range("B1").Value = 1.544050812721252
$ws.Range("C1").Value = 1.588072061538696
$ws.Range("D1").Value = 2.070609092712402
$ws.Range("E1").Value = 3.397009611129761
